# Applies the "Updated cryptos list" data refresh: new prices/volumes for
# every coin row, plus a few rows where the ranking order (and thus the
# row a coin sits on) swapped with its neighbour.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value that must stay text even when it looks like
# a number (e.g. "0.9974" or "30.591.87"). A leading apostrophe is the
# standard Excel "force text" quote-prefix, exactly like typing it by hand.
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
}

# Row 2
$ws.Cells.Item(2, 4).Value = '30.591.87'
$ws.Cells.Item(2, 5).Value = '  +1.50%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.857.62'
$ws.Cells.Item(3, 5).Value = '  +1.28%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '0.9974'
$ws.Cells.Item(4, 5).Value = '  -0.38%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '235.10'
$ws.Cells.Item(5, 5).Value = '  +1.19%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '0.9980'
$ws.Cells.Item(6, 5).Value = '  -0.31%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.4727'
$ws.Cells.Item(7, 5).Value = '  +1.42%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.2779'
$ws.Cells.Item(8, 5).Value = '  +2.41%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.06438'
$ws.Cells.Item(9, 5).Value = '  +2.82%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '18.41'
$ws.Cells.Item(10, 5).Value = '  +15.28%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Cells.Item(11, 4) '0.07459'
$ws.Cells.Item(11, 5).Value = '  +0.81%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.815.67'
$ws.Cells.Item(12, 5).Value = '  -0.98%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '88.84'
$ws.Cells.Item(13, 5).Value = '  +6.63%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '5.002'
$ws.Cells.Item(14, 5).Value = '  +1.84%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '0.6344'
$ws.Cells.Item(15, 5).Value = '  +3.14%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '301.56'
$ws.Cells.Item(16, 5).Value = '  +31.77%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '30.554.76'
$ws.Cells.Item(17, 5).Value = '  +1.63%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '0.9980'
$ws.Cells.Item(18, 5).Value = '  -0.30%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +4.34%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '0.000007416'
$ws.Cells.Item(20, 5).Value = '  +2.13%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(21, 4).Value = '2.087.83'
$ws.Cells.Item(21, 5).Value = '  +0.80%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'BinanceUSD'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Cells.Item(22, 4) '0.9976'
$ws.Cells.Item(22, 5).Value = '  -0.41%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '5.056'
$ws.Cells.Item(23, 5).Value = '  +4.43%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '6.013'
$ws.Cells.Item(24, 5).Value = '  +3.55%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '167.55'
$ws.Cells.Item(25, 5).Value = '  +1.58%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '9.122'
$ws.Cells.Item(26, 5).Value = '  -0.63%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '18.98'
$ws.Cells.Item(27, 5).Value = '  +7.46%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '1.942'
$ws.Cells.Item(28, 5).Value = '  +4.29%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '0.1039'
$ws.Cells.Item(29, 5).Value = '  +1.17%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '1.346'
$ws.Cells.Item(30, 5).Value = '  -1.83%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '4.064'
$ws.Cells.Item(31, 5).Value = '  -0.04%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '3.903'
$ws.Cells.Item(32, 5).Value = '  +3.41%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '0.04896'
$ws.Cells.Item(33, 5).Value = '  +2.50%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '1.159'
$ws.Cells.Item(34, 5).Value = '  +2.60%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '0.7125'
$ws.Cells.Item(35, 5).Value = '  +0.83%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '2.687'
$ws.Cells.Item(36, 5).Value = '  -0.95%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '0.01899'
$ws.Cells.Item(37, 5).Value = '  +1.90%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '2.683'
$ws.Cells.Item(38, 5).Value = '  +1.44%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(39, 4) '2.014'
$ws.Cells.Item(39, 5).Value = '  +4.54%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Cells.Item(40, 4) '0.8819'
$ws.Cells.Item(40, 5).Value = '  -0.83%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '106.58'
$ws.Cells.Item(41, 5).Value = '  +2.44%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '0.9977'
$ws.Cells.Item(42, 5).Value = '  -0.37%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '0.4131'
$ws.Cells.Item(43, 5).Value = '  +3.50%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '5.508'
$ws.Cells.Item(44, 5).Value = '  +0.56%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '7.290'
$ws.Cells.Item(45, 5).Value = '  +5.27%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '63.16'
$ws.Cells.Item(46, 5).Value = '  +5.80%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.1205'
$ws.Cells.Item(47, 5).Value = '  +1.69%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '34.50'
$ws.Cells.Item(48, 5).Value = '  +6.41%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '8.649'
$ws.Cells.Item(49, 5).Value = '  +1.19%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '1.370'
$ws.Cells.Item(50, 5).Value = '  +1.30%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '0.05507'
$ws.Cells.Item(51, 5).Value = '  -0.02%  '
